# Insert a new daily price record for Puerro (Azul de Maquehue, Primera)
# at Vega Modelo de Temuco as row 16, pushing the existing rows 16-142
# down to 17-143 (weekly -> new entry added at the top of the data block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (16..142) down by one row.
$ws.Rows.Item(16).EntireRow.Insert()

# Populate the newly inserted row 16 with the new record's data.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 44490
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112005
$ws.Range("G16").Value = "Puerro"
$ws.Range("H16").Value = "Azul de Maquehue"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("N16").Value = "$/docena de paquetes"
$ws.Range("O16").Value = "Provincia de Cautín"
$ws.Range("P16").Value = 583
$ws.Range("Q16").Value = 12
$ws.Range("R16").Value = "Hortaliza"
